$d = $word.ActiveDocument

# Locate the paragraph that ends the "Improvements" section (the one
# discussing the Java code / N-Grams) via a scan, so we don't depend on a
# brittle hard-coded paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Tweets could be possible cut or too long for Twitter to handle.*") {
        $target = $p
    }
}

$quote = [char]8220
$unquote = [char]8221

# Insert the three new (still-blank) paragraphs first, all derived from the
# un-bolded "target" paragraph, so none of them pick up bold formatting by
# inheritance. We apply the bold / indent tweaks afterwards.
$r1 = $target.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$p1 = $target.Next()

$r2 = $p1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p2 = $p1.Next()

$r3 = $p2.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p3 = $p2.Next()

# --- New paragraph 1: "GitHub Link: " (bold heading, first-line indent) ---
$p1.Range.Text = "GitHub Link: "
$p1.Range.Font.Name = "Times New Roman"
$p1.Range.Font.Size = 12
$p1.Range.Font.Bold = 1
$p1.Format.FirstLineIndent = 36
$p1.Format.Alignment = 3
$p1.Format.LineSpacingRule = 2
$p1.Format.LineSpacing = 24

# --- New paragraph 2: "All files are available..." (regular text) ---
$p2.Range.Text = "All files are available for download on our GitHub under the name " + $quote + "Shazbro" + $unquote + " in the " + $quote + "Trump Twitter Generator" + $unquote + " project section."
$p2.Range.Font.Name = "Times New Roman"
$p2.Range.Font.Size = 12
$p2.Format.FirstLineIndent = 36
$p2.Format.Alignment = 3
$p2.Format.LineSpacingRule = 2
$p2.Format.LineSpacing = 24

# --- New paragraph 3: "URL: ..." (no first-line indent) ---
$p3.Range.Text = "URL: https://github.com/Shazbro/Trump-Twitter-Generator"
$p3.Range.Font.Name = "Times New Roman"
$p3.Range.Font.Size = 12
$p3.Format.FirstLineIndent = 0
$p3.Format.LeftIndent = 0
$p3.Format.Alignment = 3
$p3.Format.LineSpacingRule = 2
$p3.Format.LineSpacing = 24
